$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grocery list")
$ws.Range("E2:E11").Formula = "=C2*D2"
$ws.Range("E12").Formula = "=SUM(E2:E11)"
$ws.Range("E2:E12").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$ws.Range("A2:A11").Validation.Delete()
$ws.Application.ActiveWindow.Zoom = 130
